$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target data (rows re-sorted by age, with some values changed) for rows 2-15
$data = @(
    @("Error",     120, "Córdoba",      $true,  ""),
    @("José",      65,  "Buenos Aires", $true,  "Mayor"),
    @("Carlos",    45,  "Salta",        $true,  "Adulto"),
    @("Juan",      35,  "Córdoba",      $true,  "Adulto"),
    @("Agustina",  30,  "Buenos Aires", $true,  "Jóven"),
    @("Ana",       30,  "Buenos Aires", $true,  "Jóven"),
    @("Martina",   27,  "Mendoza",      $true,  "Jóven"),
    @("Luis",      25,  "Mendoza",      $true,  "Jóven"),
    @("Tomás",     22,  "Mendoza",      $true,  "Jóven"),
    @("Nicolás",   21,  "Rosario",      $true,  "Jóven"),
    @("Valentina", 20,  "Córdoba",      $true,  "Jóven"),
    @("Lucía",     19,  "Córdoba",      $true,  "Jóven"),
    @("Federico",  18,  "Salta",        $false, "Adolescente"),
    @("Sofía",     16,  "Rosario",      $false, "Adolescente")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
